# Generate Report for Handoff
# Adds two new localization entries (34bfdbb5-... and da9e7863-...) as new
# rows in each of the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$repo = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob"
$repoZhCn = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob"
$repoDeDe = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob"
$sha = "c7a1f2e4b9d8536a0f4de1b7c9a3e5d2f6b8c0a1"

# ---------------------------------------------------------------------------
# New source files being handed off in this run.
# ---------------------------------------------------------------------------
$file1 = "34bfdbb5-28e8-43bc-b318-a409e2af5021.md"
$file1Path = "e2e\" + $file1
$file1HandoffZhCn = "34bfdbb5-28e8-43bc-b318-a409e2af5021.d652f47260306d521fd35bd955b3702df5e67f56.zh-cn.xlf"
$file1HandoffDeDe = "34bfdbb5-28e8-43bc-b318-a409e2af5021.d652f47260306d521fd35bd955b3702df5e67f56.de-de.xlf"
$file1DateDeDe = "2016-09-08 04:55:33"
$file1DateZhCn = "2016-09-08 04:55:28"

$file2 = "da9e7863-03e8-4548-a100-4732be3f6675.md"
$file2Path = "e2e\" + $file2
$file2HandoffZhCn = "da9e7863-03e8-4548-a100-4732be3f6675.f871553eef17a673d4c8599ddd41beb2d8f8319f.zh-cn.xlf"
$file2HandoffDeDe = "da9e7863-03e8-4548-a100-4732be3f6675.f871553eef17a673d4c8599ddd41beb2d8f8319f.de-de.xlf"
$file2DateDeDe = "2016-09-08 04:55:33"
$file2DateZhCn = "2016-09-08 04:55:28"

$status = "Ready for handoff"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ov = $wb.Worksheets.Item("Overview")
$ovTable = $ov.ListObjects.Item(1)
$ovTable.ListRows.Add() | Out-Null
$ovTable.ListRows.Add() | Out-Null

# Row 6 - file1
$ov.Range("A6").Value = $file1
$ov.Range("C6").Value = ".md"
$ov.Range("E6").Value = $status
$ov.Range("F6").Value = $status
$ov.Range("G6").Value = $file1DateDeDe
$ov.Range("G6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ov.Hyperlinks.Add($ov.Range("B6"), "$repo/$sha/$file1Path", "", "", $file1Path)

# Row 7 - file2
$ov.Range("A7").Value = $file2
$ov.Range("C7").Value = ".md"
$ov.Range("E7").Value = $status
$ov.Range("F7").Value = $status
$ov.Range("G7").Value = $file2DateDeDe
$ov.Range("G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ov.Hyperlinks.Add($ov.Range("B7"), "$repo/$sha/$file2Path", "", "", $file2Path)

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$zh = $wb.Worksheets.Item("zh-cn")
$zhTable = $zh.ListObjects.Item(1)
$zhTable.ListRows.Add() | Out-Null
$zhTable.ListRows.Add() | Out-Null

# Row 6 - file1
$zh.Range("B6").Value = ".md"
$zh.Range("C6").Value = $status
$zh.Range("D6").Value = "e2e"
$zh.Range("E6").Value = "ht"
$zh.Range("F6").Value = "False"
$zh.Range("G6").Value = $file1HandoffZhCn
$zh.Range("H6").Value = $file1DateZhCn
$zh.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("K6").Value = "0001-01-01 00:00:00"
$zh.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("M6").Value = "True"
$zh.Range("O6").Value = "False"
$zh.Hyperlinks.Add($zh.Range("A6"), "$repo/$sha/$file1Path", "", "", $file1)

# Row 7 - file2
$zh.Range("B7").Value = ".md"
$zh.Range("C7").Value = $status
$zh.Range("D7").Value = "e2e"
$zh.Range("E7").Value = "ht"
$zh.Range("F7").Value = "False"
$zh.Range("G7").Value = $file2HandoffZhCn
$zh.Range("H7").Value = $file2DateZhCn
$zh.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("K7").Value = "0001-01-01 00:00:00"
$zh.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("M7").Value = "True"
$zh.Range("O7").Value = "False"
$zh.Hyperlinks.Add($zh.Range("A7"), "$repo/$sha/$file2Path", "", "", $file2)

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$de = $wb.Worksheets.Item("de-de")
$deTable = $de.ListObjects.Item(1)
$deTable.ListRows.Add() | Out-Null
$deTable.ListRows.Add() | Out-Null

# Row 6 - file1
$de.Range("B6").Value = ".md"
$de.Range("C6").Value = $status
$de.Range("D6").Value = "e2e"
$de.Range("E6").Value = "ht"
$de.Range("F6").Value = "False"
$de.Range("G6").Value = $file1HandoffDeDe
$de.Range("H6").Value = $file1DateDeDe
$de.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("K6").Value = "0001-01-01 00:00:00"
$de.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("M6").Value = "True"
$de.Range("O6").Value = "False"
$de.Hyperlinks.Add($de.Range("A6"), "$repo/$sha/$file1Path", "", "", $file1)

# Row 7 - file2
$de.Range("B7").Value = ".md"
$de.Range("C7").Value = $status
$de.Range("D7").Value = "e2e"
$de.Range("E7").Value = "ht"
$de.Range("F7").Value = "False"
$de.Range("G7").Value = $file2HandoffDeDe
$de.Range("H7").Value = $file2DateDeDe
$de.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("K7").Value = "0001-01-01 00:00:00"
$de.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("M7").Value = "True"
$de.Range("O7").Value = "False"
$de.Hyperlinks.Add($de.Range("A7"), "$repo/$sha/$file2Path", "", "", $file2)

Write-Output "Report generated for handoff: added $file1 and $file2"
